# The document's (only) footer started with a small "-<page>-" page-number
# construct: a literal "-" run, then a content control (structured document
# tag) wrapping a PAGE field ("-2-" in total), then the paragraph mark.
# The edit removes that whole construct, leaving the first footer paragraph
# empty (the second, already-empty footer paragraph is left untouched).

$d = $word.ActiveDocument
$footer = $d.Sections(1).Footers(1)

# Remove the page-number content control together with its field content
# (the "PAGE \* MERGEFORMAT" field and its cached result).
$cc = $footer.Range.ContentControls(1)
$cc.Delete($true)

# Remove the leftover literal "-" run that preceded the content control,
# without touching the paragraph mark, so the paragraph itself survives
# (now empty) exactly like in the target document.
$p1 = $footer.Range.Paragraphs(1)
$r = $p1.Range
$r.End = $r.End - 1
$r.Text = ""
